$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'275.64"
$ws.Range("D4").Formula = "'6.397"
$ws.Range("D5").Formula = "'0.06272"
$ws.Range("D6").Formula = "'3.657"
$ws.Range("D7").Formula = "'6.665"
$ws.Range("D8").Formula = "'1.355"
$ws.Range("D9").Formula = "'0.8306"
$ws.Range("D10").Formula = "'0.01376"
$ws.Range("D11").Formula = "'0.1622"
$ws.Range("D12").Formula = "'0.08303"
$ws.Range("D13").Formula = "'0.03445"
$ws.Range("D14").Formula = "'0.03078"
$ws.Range("D15").Formula = "'0.09311"
$ws.Range("D16").Formula = "'3.854"
$ws.Range("D17").Formula = "'0.001644"
$ws.Range("D18").Formula = "'0.04775"
$ws.Range("D19").Formula = "'0.006310"
$ws.Range("D20").Formula = "'0.005690"
$ws.Range("D23").Formula = "'3.717"
$ws.Range("D27").Formula = "'0.0002679"
$ws.Range("D40").Formula = "'0.04697"
$ws.Range("D41").Formula = "'0.007054"
$ws.Range("D42").Formula = "'0.1163"
$ws.Range("D43").Formula = "'0.003349"
$ws.Range("D45").Formula = "'0.00006248"
$ws.Range("D48").Formula = "'0.7696"
$ws.Range("D49").Formula = "'0.03997"
$ws.Range("D50").Formula = "'0.00002299"
$ws.Range("D51").Formula = "'0.01239"
